$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 32) with the latest portfolio values.
# The date column must stay plain text (not get auto-converted to a
# date serial number by Excel's input parser), so we briefly force a
# text number format while assigning the value, then clear the
# formatting again so the new cell ends up unstyled like its peers.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2025-09-16"
$ws.Range("A32").ClearFormats()

$ws.Range("B32").Value = 58.97999954223633
$ws.Range("C32").Value = 713.25
$ws.Range("D32").Value = 327.2999877929688
